$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new risk row (row 8). Strings are written in the same order they
# first appear in the target sharedStrings table so the new <si> entries
# land at indices 28-31 in that order (B8, G8, H8, A8).
$ws.Range("B8").Value = "As there is no password or secure login, anyone who knows your username can log into your account. Should not be used to store sensitive information."
$ws.Range("G8").Value = "Make sure people know not to use it to store sensitive information"
$ws.Range("H8").Value = "Take down sensitive infromation immediately and possibly inform authorities/relevant parties."
$ws.Range("A8").Value = "Someone getting sensitive information from your account"

$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = "Developer/user"

# Update the sheet view: scroll so column D is the left-most visible column,
# set normal zoom to 100% (and drop the custom zoomScale), and move the
# active selection to F14.
$ws.Select()
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F14").Select()
